$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.333.37'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -2.85%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.941.27'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -2.90%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.26'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -2.36%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7234'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  -7.57%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3361'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  -4.59%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '28.70'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07439'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +5.67%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8163'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  -5.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08146'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  -0.70%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.938.01'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  -3.01%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.517'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  -1.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '95.29'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  -5.70%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.87'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -4.83%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.357.45'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008327'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  +4.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '254.41'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  -7.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.874'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  -1.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.194.81'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  -3.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +0.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.954'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  -2.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.866'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  -2.16%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.62'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -2.43%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.435'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  +2.47%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.43'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  -2.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1322'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  -10.69%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.573'
$ws.Range('D30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.342'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -1.55%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.487'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  -2.85%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.239'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -5.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05255'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  +0.75%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.273'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  +3.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7562'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -2.87%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.740'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  -2.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01990'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -0.91%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.836'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -3.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '81.61'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  +2.33%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.555'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  -3.17%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4580'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  -2.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.025'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  -6.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8488'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -0.25%  '

$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.01'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  -3.29%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.883'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  -1.15%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.432'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  -3.86%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.11'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +0.66%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4200'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  -3.15%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.507'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  -0.57%  '
